# Update of all scripts and data
# - remove three incomplete-catch rows (missing Numb) from the "2-RAP" gear block
# - fill in the weight/number for Pecten jacobaeus (2-RAP) that was missing them
# - recompute the raising factor (RF) for the benthos/discard rows that previously had none
# - correct Numb for the "discard"/concretion-type rows to -1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for "Eledone moschata" (18), "Galeodea echinophora" (19) and
# "Squilla mantis" (28). Work from the bottom up so row numbers for the
# earlier deletions stay valid.
$ws.Rows("28:28").Delete()
$ws.Rows("19:19").Delete()
$ws.Rows("18:18").Delete()

# After the deletion, "Pecten jacobaeus" (2-RAP) lands on row 23 and is missing
# its weight/number values.
$ws.Range("G23").Value = 0.057
$ws.Range("H23").Value = 1

# Rows 26-44 (2-RAP benthos / discard / debris rows) get a raising factor (RF).
$ws.Range("I26:I44").Value = 34.02508771929825

# Fix the catch number for these rows to -1.
$ws.Range("H29").Value = -1
$ws.Range("H30").Value = -1
$ws.Range("H39").Value = -1
$ws.Range("H44").Value = -1
